# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Sat Dec 30 22:17:26 UTC 2023 with GitHub Actions".
#
# All data cells on the sheet are plain text (coin name / link / price /
# 1h-volume), even though several of the new "Price" values look like plain
# decimal numbers (e.g. "316.60"). Setting such a string straight into
# `.Value` would make Excel auto-convert it to a numeric cell (and drop the
# trailing zero), so for column D we briefly mark the cell as Text (`@`)
# before writing the value, then restore the original "Normal" style so the
# cell formatting/look is unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.460.76'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.06%  '
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.300.89'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.90%  '
# Row 4
$ws.Range("E4").Value = '  -0.10%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.60'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.86%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.87'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.93%  '
# Row 7
$ws.Range("E7").Value = '  +0.78%  '
# Row 8
$ws.Range("E8").Value = '  -0.02%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.606'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.45%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.70'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.33%  '
# Row 11
$ws.Range("E11").Value = '  -0.40%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.43'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.36%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.107'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.67%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.963'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.04%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.27'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.75%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.647.82'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.76%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.295.91'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.19%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.329.07'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.89%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.42'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.36%  '
# Row 20
$ws.Range("E20").Value = '  +1.31%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.45'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.80%  '
# Row 22
$ws.Range("B22").Value = 'InternetComputer(DFINITY)'
$ws.Range("C22").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.81'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +27.26%  '
# Row 23
$ws.Range("B23").Value = 'PancakeSwap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.55'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.79%  '
# Row 24
$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '273.93'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +7.06%  '
# Row 25
$ws.Range("E25").Value = '  -1.01%  '
# Row 26
$ws.Range("E26").Value = '  -0.25%  '
# Row 27
$ws.Range("E27").Value = '  -0.66%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.37'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.73%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.77'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.56%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.53'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.61%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '166.27'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.07%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0877'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.17%  '
# Row 33
$ws.Range("E33").Value = '  +1.90%  '
# Row 34
$ws.Range("E34").Value = '  +4.59%  '
# Row 35
$ws.Range("E35").Value = '  +1.12%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.62'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -10.19%  '
# Row 37
$ws.Range("E37").Value = '  +1.70%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0365'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.88%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.74'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.42%  '
# Row 40
$ws.Range("E40").Value = '  +0.95%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.50'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.70%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '70.18'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.34%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '94.60'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.36%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.226'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.11%  '
# Row 45
$ws.Range("E45").Value = '  -0.26%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.08'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.87%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '80.89'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +9.17%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '112.74'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.25%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.97'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.18%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.22'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.21%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.587.64'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.68%  '
